$d = $word.ActiveDocument

# --- 1) Merge the two split runs holding the binary string into a single
#        run with the same (concatenated) text. ---
$null = $d.Content.Find.Execute(
    "0001100010000000000001010000", $true, $false, $false, $false, $false,
    $true, 1, $false, "0001100010000000000001010000", 2)

# --- 2) "Reg A :3  Reg B : 2" -> "Reg A :-3  Reg B : 2"
#        Insert a "-" right before the "3" (i.e. turn the answer negative).
#        Do the text substitution scoped to exactly the run that holds
#        "3  Reg" so the preceding run (tab + "Reg A :") is left intact. ---
$null = $d.Content.Find.Execute(
    "3  Reg", $true, $false, $false, $false, $false,
    $true, 1, $false, "-3  Reg", 2)

# --- 3) Move the "_GoBack" bookmark from the end of the document to right
#        after the edit we just made (where the cursor was left). ---
$rng = $d.Content
$null = $rng.Find.Execute("Reg A :")
$rng.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rng)
